$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$rows = @(
    @("cs", "lab.vape.edit.title", "Editace vapu"),
    @("cs", "lab.vape.edit.subtitle", "Každý se někdy překlepne, zde je možné upravit vape."),
    @("cs", "lab.vape.link.button", "Zpět"),
    @("cs", "lab.vape.update.submit", "Aktualizovat"),
    @("cs", "lab.vape.update.success", "Vape byl aktualizován."),
    @("cs", "lab.vape.index.title", "Náhled vapu"),
    @("cs", "lab.vape.button.edit", "Editovat"),
    @("cs", "lab.vape.index.preview.subtitle", "Správa vybraného vapu"),
    @("cs", "lab.vape.index.preview.title", "Náhled vapu")
)

$startRow = 520
$endRow = $startRow + $rows.Count - 1

# Copy the formatting of the last existing data row onto the new rows
# before filling in values, so the new cells pick up the same style
# index as the rest of the table instead of the engine's default.
$ws.Range("A519:C519").Copy()
$ws.Range("A" + $startRow + ":C" + $endRow).PasteSpecial(-4122)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$ws.Range("B524").Select()
